$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.291.60'
$ws.Range("E2").Value = '  -3.02%  '
$ws.Range("D3").Value = '1.854.96'
$ws.Range("E3").Value = '  -3.68%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '323.49'
$ws.Range("E5").Value = '  -1.74%  '
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '0.4522'
$ws.Range("E7").Value = '  -4.35%  '
$ws.Range("D8").Value = '0.3859'
$ws.Range("E8").Value = '  -4.83%  '
$ws.Range("D9").Value = '48.69'
$ws.Range("E9").Value = '  -8.10%  '
$ws.Range("D10").Value = '0.07928'
$ws.Range("E10").Value = '  -6.10%  '
$ws.Range("D11").Value = '1.013'
$ws.Range("E11").Value = '  -3.22%  '
$ws.Range("D12").Value = '21.33'
$ws.Range("E12").Value = '  -4.01%  '
$ws.Range("D13").Value = '1.848.41'
$ws.Range("E13").Value = '  -3.43%  '
$ws.Range("D14").Value = '5.902'
$ws.Range("E14").Value = '  -3.10%  '
$ws.Range("D15").Value = '7.127'
$ws.Range("E15").Value = '  -5.06%  '
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '85.82'
$ws.Range("E17").Value = '  -5.02%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '0.06570'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.00001025'
$ws.Range("E19").Value = '  -3.63%  '
$ws.Range("E20").Value = '  -6.15%  '
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = '5.505'
$ws.Range("E22").Value = '  -4.52%  '
$ws.Range("D23").Value = '27.289.34'
$ws.Range("E23").Value = '  -3.06%  '
$ws.Range("D24").Value = '10.87'
$ws.Range("E24").Value = '  -4.77%  '
$ws.Range("D25").Value = '2.288'
$ws.Range("E25").Value = '  +0.56%  '
$ws.Range("D26").Value = '2.077.59'
$ws.Range("E26").Value = '  -3.04%  '
$ws.Range("D27").Value = '153.63'
$ws.Range("E27").Value = '  -0.60%  '
$ws.Range("D28").Value = '19.88'
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("D29").Value = '2.062'
$ws.Range("E29").Value = '  -4.14%  '
$ws.Range("D30").Value = '5.454'
$ws.Range("E30").Value = '  -5.03%  '
$ws.Range("D31").Value = '121.18'
$ws.Range("E31").Value = '  -2.15%  '
$ws.Range("D32").Value = '0.09298'
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("D33").Value = '0.9321'
$ws.Range("E33").Value = '  -4.63%  '
$ws.Range("D34").Value = '1.460'
$ws.Range("E34").Value = '  +1.59%  '
$ws.Range("D35").Value = '3.585'
$ws.Range("E35").Value = '  -1.49%  '
$ws.Range("D36").Value = '5.265'
$ws.Range("E36").Value = '  -5.23%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02222'
$ws.Range("E37").Value = '  -3.89%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.05993'
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '1.220'
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("D40").Value = '8.088'
$ws.Range("E40").Value = '  -10.91%  '
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").Value = '0.5902'
$ws.Range("E42").Value = '  -4.33%  '
$ws.Range("D43").Value = '0.1883'
$ws.Range("E43").Value = '  -1.00%  '
$ws.Range("D44").Value = '10.12'
$ws.Range("E44").Value = '  -8.43%  '
$ws.Range("E45").Value = '  -2.55%  '
$ws.Range("D46").Value = '0.5637'
$ws.Range("E46").Value = '  -4.28%  '
$ws.Range("D47").Value = '12.05'
$ws.Range("E47").Value = '  -6.32%  '
$ws.Range("D48").Value = '3.373'
$ws.Range("E48").Value = '  -2.72%  '
$ws.Range("D49").Value = '1.913'
$ws.Range("E49").Value = '  -5.87%  '
$ws.Range("D50").Value = '0.06739'
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("D51").Value = '108.61'
$ws.Range("E51").Value = '  -1.06%  '
